$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.788.22'
$ws.Cells.Item(2, 5).Value = '  +1.41%  '

$ws.Cells.Item(3, 4).Value = '3.159.55'
$ws.Cells.Item(3, 5).Value = '  +2.62%  '

$ws.Cells.Item(4, 5).Value = '  -0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '573.37'
$ws.Cells.Item(5, 5).Value = '  +2.63%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '151.22'
$ws.Cells.Item(6, 5).Value = '  +5.61%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  -0.16%  '

$ws.Cells.Item(8, 4).Value = '3.155.01'
$ws.Cells.Item(8, 5).Value = '  +2.61%  '

$ws.Cells.Item(9, 5).Value = '  +4.96%  '

$ws.Cells.Item(10, 5).Value = '  +6.90%  '

$ws.Cells.Item(11, 5).Value = '  +0.85%  '

$ws.Cells.Item(12, 5).Value = '  +7.71%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000256'
$ws.Cells.Item(13, 5).Value = '  +12.94%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '37.97'
$ws.Cells.Item(14, 5).Value = '  +8.78%  '

$ws.Cells.Item(15, 4).Value = '3.673.27'
$ws.Cells.Item(15, 5).Value = '  +2.31%  '

$ws.Cells.Item(16, 4).Value = '64.924.01'
$ws.Cells.Item(16, 5).Value = '  +1.44%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '7.25'
$ws.Cells.Item(17, 5).Value = '  +7.85%  '

$ws.Cells.Item(18, 4).Value = '3.156.82'
$ws.Cells.Item(18, 5).Value = '  +2.49%  '

$ws.Cells.Item(19, 5).Value = '  +0.64%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '515.28'
$ws.Cells.Item(20, 5).Value = '  +7.53%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '15.02'
$ws.Cells.Item(21, 5).Value = '  +8.04%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.738'
$ws.Cells.Item(22, 5).Value = '  +10.24%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '15.15'
$ws.Cells.Item(23, 5).Value = '  +8.05%  '

$ws.Cells.Item(24, 5).Value = '  +4.46%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '85.27'
$ws.Cells.Item(25, 5).Value = '  +5.28%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 5).Value = '  +0.07%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.94'
$ws.Cells.Item(27, 5).Value = '  +5.32%  '

$ws.Cells.Item(28, 5).Value = '  +10.12%  '

$ws.Cells.Item(29, 5).Value = '  +6.65%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '28.04'
$ws.Cells.Item(30, 5).Value = '  +6.78%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.00'
$ws.Cells.Item(31, 5).Value = '  -0.21%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.18'
$ws.Cells.Item(32, 5).Value = '  +3.74%  '

$ws.Cells.Item(33, 5).Value = '  +7.86%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '6.13'
$ws.Cells.Item(34, 5).Value = '  +10.37%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '6.59'
$ws.Cells.Item(35, 5).Value = '  +6.83%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '55.73'
$ws.Cells.Item(36, 5).Value = '  +0.00%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '486.68'
$ws.Cells.Item(37, 5).Value = '  +8.22%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0865'
$ws.Cells.Item(38, 5).Value = '  +5.96%  '

$ws.Cells.Item(39, 5).Value = '  +4.31%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.02'
$ws.Cells.Item(40, 5).Value = '  +3.36%  '

$ws.Cells.Item(41, 4).Value = '3.118.61'
$ws.Cells.Item(41, 5).Value = '  +5.00%  '

$ws.Cells.Item(42, 5).Value = '  +5.55%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.121'
$ws.Cells.Item(43, 5).Value = '  +5.82%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.296'
$ws.Cells.Item(44, 5).Value = '  +14.33%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.46'
$ws.Cells.Item(45, 5).Value = '  +16.55%  '

$ws.Cells.Item(46, 5).Value = '  +5.35%  '

$ws.Cells.Item(47, 4).Value = '0.0₃0581'
$ws.Cells.Item(47, 5).Value = '  +13.94%  '

$ws.Cells.Item(49, 5).Value = '  +3.98%  '

$ws.Cells.Item(50, 5).Value = '  +11.02%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '118.83'
$ws.Cells.Item(51, 5).Value = '  -0.55%  '

